# "puse los usuarios y contraseñas correctas"
# Replace the plaintext-password table with the correct SHA-256 password
# hashes for every user, drop the now-unused plaintext "Contraseña" column,
# and tidy up the layout (column B width / selection) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Put the correct password hashes in column B for each user (rows 2-5).
# Rows 4 and 5 previously had no hash at all - they get one now.
$ws.Range("B2").Value = "62E6BEA9285CD5040EBEB8349FD37884C58FC489083B67AB58F29B142726502B"
$ws.Range("B3").Value = "72FB1E6F1436D5A08558E9797A2B048E48468CF7596214E35C8175FD07528F4F"
$ws.Range("B4").Value = "D749DD4F4D1390494C0DA15BB5BD1CEE811DA43D6C9D45EA30AA89D95EACC5A5"
$ws.Range("B5").Value = "CDF5A0B901FD777EB7A4798AD0872C67546942FA2FEDDD70DE85508162943F31"

# The plaintext "Contraseña" column (E1 header, G2:G3 values) is no longer
# needed now that real hashes are stored - clear it out, including the
# leftover table border around the old E1 header cell.
$ws.Range("E1").ClearContents()
$ws.Range("E1").Borders.LineStyle = -4142
$ws.Range("G1").ClearContents()
$ws.Range("G2").ClearContents()
$ws.Range("G3").ClearContents()

# B2 picks up the plain "Normal" look (no leftover custom border style)
# now that it just holds a hash like every other row.
$ws.Range("B2").Style = "Normal"

# Column B needs to be a bit wider (no longer auto "best fit") to
# comfortably show the new hash values.
$ws.Range("B:B").ColumnWidth = 72.6

# Leave the selection where the user ended up after entering the data.
$ws.Range("B9").Select()
